# Apply updated crypto price/volume data (GitHub Actions scheduled refresh)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($row, $col, $text) {
    # Prefix with a leading apostrophe so Excel stores the value as literal
    # text (matching the workbook's existing inline-string cells) instead of
    # re-interpreting number-like strings (e.g. "66.658.29", "3.90", "0.0415").
    $ws.Cells.Item($row, $col).Value = "`'" + $text
}

Set-TextCell 2 4 "66.658.29"
Set-TextCell 2 5 "  +8.53%  "
Set-TextCell 3 4 "3.494.01"
Set-TextCell 3 5 "  +11.73%  "
Set-TextCell 4 5 "  -0.05%  "
Set-TextCell 5 4 "188.08"
Set-TextCell 5 5 "  +12.19%  "
Set-TextCell 6 4 "548.57"
Set-TextCell 6 5 "  +7.89%  "
Set-TextCell 7 4 "3.484.90"
Set-TextCell 7 5 "  +11.71%  "
Set-TextCell 8 4 "0.605"
Set-TextCell 8 5 "  +4.73%  "
Set-TextCell 9 5 "  -0.07%  "
Set-TextCell 10 5 "  +8.39%  "
Set-TextCell 11 5 "  +19.11%  "
Set-TextCell 12 4 "55.10"
Set-TextCell 12 5 "  +7.76%  "
Set-TextCell 13 5 "  +9.17%  "
Set-TextCell 14 5 "  +8.44%  "
Set-TextCell 15 4 "4.058.88"
Set-TextCell 15 5 "  +11.69%  "
Set-TextCell 16 4 "3.492.17"
Set-TextCell 16 5 "  +11.51%  "
Set-TextCell 17 5 "  +8.13%  "
Set-TextCell 18 4 "66.604.17"
Set-TextCell 18 5 "  +8.50%  "
Set-TextCell 19 4 "18.19"
Set-TextCell 19 5 "  +9.06%  "
Set-TextCell 20 5 "  +11.59%  "
Set-TextCell 21 5 "  +6.22%  "
Set-TextCell 22 4 "413.08"
Set-TextCell 22 5 "  +16.37%  "
Set-TextCell 23 4 "85.10"
Set-TextCell 23 5 "  +8.15%  "
Set-TextCell 24 4 "3.90"
Set-TextCell 24 5 "  +8.28%  "
Set-TextCell 25 5 "  +11.88%  "
Set-TextCell 26 4 "11.12"
Set-TextCell 26 5 "  +4.58%  "
Set-TextCell 27 5 "  +15.86%  "
Set-TextCell 28 4 "6.11"
Set-TextCell 28 5 "  +0.44%  "
Set-TextCell 29 4 "11.84"
Set-TextCell 29 5 "  +9.27%  "
Set-TextCell 30 4 "8.82"
Set-TextCell 30 5 "  +11.77%  "
Set-TextCell 31 4 "30.21"
Set-TextCell 31 5 "  +9.88%  "
Set-TextCell 32 4 "654.10"
Set-TextCell 32 5 "  +3.13%  "
Set-TextCell 33 5 "  +7.77%  "
Set-TextCell 34 4 "11.69"
Set-TextCell 34 5 "  +7.02%  "
Set-TextCell 35 5 "  +10.26%  "
Set-TextCell 36 4 "59.62"
Set-TextCell 36 5 "  +7.45%  "
Set-TextCell 37 4 "38.63"
Set-TextCell 37 5 "  +9.85%  "
Set-TextCell 38 4 "0.0₃0811"
Set-TextCell 38 5 "  +20.31%  "
Set-TextCell 39 5 "  -0.05%  "
Set-TextCell 40 5 "  +8.19%  "
Set-TextCell 41 5 "  +13.75%  "
Set-TextCell 42 5 "  +20.56%  "
Set-TextCell 43 5 "  +0.08%  "
Set-TextCell 44 4 "2.988.68"
Set-TextCell 44 5 "  +7.60%  "
Set-TextCell 45 2 "ThetaToken"
Set-TextCell 45 3 "https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta"
Set-TextCell 45 4 "2.90"
Set-TextCell 45 5 "  +17.78%  "
Set-TextCell 46 2 "Fetch.AI"
Set-TextCell 46 3 "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
Set-TextCell 46 4 "2.63"
Set-TextCell 46 5 "  +9.10%  "
Set-TextCell 47 4 "3.28"
Set-TextCell 47 5 "  +13.61%  "
Set-TextCell 48 4 "0.0415"
Set-TextCell 48 5 "  +11.02%  "
Set-TextCell 49 4 "2.71"
Set-TextCell 49 5 "  +5.24%  "
Set-TextCell 50 4 "8.91"
Set-TextCell 50 5 "  +21.67%  "
Set-TextCell 51 5 "  +8.60%  "
